{"js": "// Add four new rows to the end of the \"Test Log\" table describing the\n// \"blue background\" checks for the easy / medium / hard / insane game\n// options, mirroring the existing table's row layout (Action, Expected\n// Result, Pass/Fail, Reason for Failure).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in the document body.\");\n}\n\nconst table = tables.items[0];\n\nconst newRows = [\n  [\n    \"The water has blue background in relevant areas as regards the easy option of the game \",\n    \"No\",\n    \"Pass\",\n    \"None\"\n  ],\n  [\n    \"The water has blue background in relevant areas as regards the medium option of the game\",\n    \"Yes\",\n    \"Pass\",\n    \"None\"\n  ],\n  [\n    \"The water has blue background in relevant areas as regards the hard option of the game\",\n    \"Yes\",\n    \"Pass\",\n    \"None\"\n  ],\n  [\n    \"The water has blue background in relevant areas as regards the insane option of the game\",\n    \"Yes\",\n    \"Fail\",\n    \"Feature not operational\"\n  ]\n];\n\ntable.addRows(\"End\", newRows.length, newRows);\nawait context.sync();\n", "ps1": "# Add four new rows to the end of the \"Test Log\" table describing the\n# \"blue background\" checks for the easy / medium / hard / insane game\n# options, mirroring the existing table's row layout (Action, Expected\n# Result, Pass/Fail, Reason for Failure).\n$doc = $word.ActiveDocument\n$table = $doc.Tables.Item(1)\n\n$newRows = @(\n    @(\"The water has blue background in relevant areas as regards the easy option of the game \", \"No\", \"Pass\", \"None\"),\n    @(\"The water has blue background in relevant areas as regards the medium option of the game\", \"Yes\", \"Pass\", \"None\"),\n    @(\"The water has blue background in relevant areas as regards the hard option of the game\", \"Yes\", \"Pass\", \"None\"),\n    @(\"The water has blue background in relevant areas as regards the insane option of the game\", \"Yes\", \"Fail\", \"Feature not operational\")\n)\n\nforeach ($rowValues in $newRows) {\n    $table.Rows.Add() | Out-Null\n    $rowIdx = $table.Rows.Count\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $table.Cell($rowIdx, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
